# Update redis & clickhouse connection examples
#
# - extend the "Desc" cell explaining the Type column with a 4th option
#   (clickhouse)
# - label the (previously blank) D10 cell under the "IP" header
# - append a new data row (14) describing a ClickhouseLogDb_1 connection,
#   mirroring the existing RedisGameDb_1 row's layout/formatting
# - move the active selection to D18 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Desc row (row 10) -----------------------------------------------------
$ws.Cells.Item(10, 3).Value = "1: mysql  2:mongodb 3: redis 4: clickhouse"
$ws.Cells.Item(10, 4).Value = "IP"

# --- New row 14: Clickhouse log DB entry -----------------------------------
# Clone row 13's formatting (cell styles/borders/number formats) first so the
# new row matches the rest of the table, then fill in the values.
$ws.Range("A13:I13").Copy()
$ws.Range("A14:I14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Cells.Item(14, 1).Value = "ClickhouseLogDb_1"
$ws.Cells.Item(14, 2).Value = "4"
$ws.Cells.Item(14, 3).Value = "4"
$ws.Cells.Item(14, 4).Value = "127.0.0.1"
$ws.Cells.Item(14, 5).Value = "1.14.123.62"
$ws.Cells.Item(14, 6).Value = 10431
$ws.Cells.Item(14, 7).Value = "default"
$ws.Cells.Item(14, 8).Value = "pwnsky_squick"
$ws.Cells.Item(14, 9).Value = "squick"

# --- View state --------------------------------------------------------
$ws.Range("D18").Select() | Out-Null
